$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 4.6
$ws.Range("N2").Value = 4.6
$ws.Range("Q2").Value = 1.76
$ws.Range("S2").Value = 2.92
$ws.Range("T2").Value = 1.73
$ws.Range("U2").Value = 2.28
$ws.Range("V2").Value = 1.27
$ws.Range("W2").Value = 2.14
$ws.Range("AM2").Value = 85
$ws.Range("AO2").Value = 46
$ws.Range("F3").Value = 3.45
$ws.Range("O4").Value = 1.27
$ws.Range("P4").Value = 1.95
$ws.Range("T4").Value = 1.04
$ws.Range("P5").Value = 1.33
$ws.Range("Q5").Value = 1.01
$ws.Range("S5").Value = 7.6
$ws.Range("M6").Value = 1.06
$ws.Range("O6").Value = 1.06
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 6
$ws.Range("I7").Value = 1.66
$ws.Range("J7").Value = 4.7
$ws.Range("K7").Value = 5.7
$ws.Range("P7").Value = 2.68
$ws.Range("Q7").Value = 1.47
$ws.Range("R7").Value = 1.68
$ws.Range("S7").Value = 2.18
$ws.Range("T7").Value = 1.6
$ws.Range("U7").Value = 2.32
$ws.Range("V7").Value = 2.52
$ws.Range("H8").Value = 1.52
$ws.Range("L8").Value = 1.42
$ws.Range("N8").Value = 3.6
$ws.Range("O8").Value = 1.33
$ws.Range("R8").Value = 1.32
$ws.Range("S8").Value = 3.55
$ws.Range("T8").Value = 2.08
$ws.Range("U8").Value = 1.76
$ws.Range("AD8").Value = 10
$ws.Range("F9").Value = 1.63
$ws.Range("H9").Value = 6.6
$ws.Range("M9").Value = 1.08
$ws.Range("R9").Value = 1.22
$ws.Range("T9").Value = 2
$ws.Range("V9").Value = 1.15
$ws.Range("X9").Value = 1000
$ws.Range("G10").Value = 2.96
$ws.Range("P10").Value = 1.94
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = 1.35
$ws.Range("W10").Value = 1.51
$ws.Range("AC10").Value = 7.4
$ws.Range("AM10").Value = 90
$ws.Range("Q11").Value = 1.97
$ws.Range("AO11").Value = 9.199999999999999
$ws.Range("W12").Value = 1.79
$ws.Range("AB12").Value = 5.9
$ws.Range("AD12").Value = 22
$ws.Range("AG12").Value = 13.5
$ws.Range("AJ12").Value = 32
$ws.Range("AL12").Value = 95
$ws.Range("AN12").Value = 42
$ws.Range("AO12").Value = 190
$ws.Range("J13").Value = 3.55
$ws.Range("L13").Value = 1.46
$ws.Range("M13").Value = 1.07
$ws.Range("N13").Value = 3.2
$ws.Range("O13").Value = 1.39
$ws.Range("Q13").Value = 1.92
$ws.Range("R13").Value = 1.26
$ws.Range("S13").Value = 3.6
$ws.Range("T13").Value = 1.01
$ws.Range("U13").Value = 1.01
$ws.Range("V13").Value = 1.22
$ws.Range("W13").Value = 2.08
$ws.Range("X13").Value = 13
$ws.Range("Y13").Value = 16.5
$ws.Range("Z13").Value = 38
$ws.Range("AA13").Value = 150
$ws.Range("AB13").Value = 8.199999999999999
$ws.Range("AC13").Value = 8.4
$ws.Range("AD13").Value = 24
$ws.Range("AE13").Value = 80
$ws.Range("AF13").Value = 11
$ws.Range("AG13").Value = 10.5
$ws.Range("AH13").Value = 23
$ws.Range("AI13").Value = 90
$ws.Range("AJ13").Value = 22
$ws.Range("AK13").Value = 22
$ws.Range("AL13").Value = 44
$ws.Range("AM13").Value = 170
$ws.Range("AN13").Value = 110
